# This script refreshes the "cryptos" price table with newly-fetched
# market data (price + 1h volume change) for each coin row, matching the
# values produced by the scheduled GitHub Actions data-refresh job.
#
# Cell values are applied from an ordered map of cell reference -> new
# text value. Because column D holds prices that are stored as *text*
# (e.g. "1.003", "0.3933") rather than numbers, we force Excel to keep
# them as text instead of silently auto-converting them to numeric
# values. This is done by temporarily prefixing the value with a
# leading apostrophe (Excel's "treat as text" marker) and then clearing
# the resulting cell formatting so no extra number-format/style gets
# attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "D2" = '21.985.03'
    "E2" = '  -1.78%  '
    "D3" = '1.550.61'
    "E3" = '  -0.95%  '
    "D4" = '1.003'
    "E4" = '  +0.16%  '
    "D5" = '1.001'
    "E5" = '  +0.07%  '
    "D6" = '285.97'
    "E6" = '  -0.06%  '
    "D7" = '0.3933'
    "E7" = '  +5.21%  '
    "D8" = '0.3184'
    "E8" = '  -2.51%  '
    "D9" = '42.05'
    "E9" = '  -6.98%  '
    "D10" = '0.07249'
    "E10" = '  -1.88%  '
    "D11" = '1.086'
    "E11" = '  -4.90%  '
    "D12" = '1.003'
    "E12" = '  +0.17%  '
    "D13" = '18.72'
    "E13" = '  -8.18%  '
    "D14" = '5.597'
    "E14" = '  -3.87%  '
    "D15" = '6.649'
    "E15" = '  -2.30%  '
    "B16" = 'WrappedEther'
    "C16" = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    "D16" = '1.553.02'
    "E16" = '  -1.06%  '
    "B17" = 'ShibaInu'
    "C17" = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    "D17" = '0.00001119'
    "E17" = '  +2.30%  '
    "D18" = '0.06581'
    "E18" = '  -1.80%  '
    "D19" = '84.01'
    "E19" = '  -1.79%  '
    "E20" = '  +0.04%  '
    "D21" = '6.282'
    "E21" = '  -0.79%  '
    "D22" = '15.65'
    "E22" = '  -3.50%  '
    "D23" = '11.17'
    "E23" = '  -4.18%  '
    "D24" = '21.989.48'
    "E24" = '  -1.73%  '
    "D25" = '2.349'
    "E25" = '  +2.15%  '
    "D26" = '2.419'
    "E26" = '  -4.36%  '
    "D27" = '146.93'
    "E27" = '  -2.04%  '
    "D28" = '18.56'
    "E28" = '  -4.11%  '
    "D29" = '4.835'
    "E29" = '  -1.42%  '
    "D30" = '1.725.74'
    "E30" = '  -1.06%  '
    "D31" = '118.66'
    "E31" = '  -3.44%  '
    "D32" = '1.060'
    "E32" = '  +0.84%  '
    "D33" = '5.630'
    "D34" = '0.08291'
    "E34" = '  +1.13%  '
    "D35" = '9.177'
    "E35" = '  -3.63%  '
    "D36" = '1.584'
    "E36" = '  -18.31%  '
    "D37" = '0.06156'
    "E37" = '  -2.21%  '
    "D38" = '0.02251'
    "E38" = '  -5.57%  '
    "D39" = '5.068'
    "E39" = '  -3.17%  '
    "D40" = '0.2062'
    "E40" = '  -5.38%  '
    "E41" = '  -6.37%  '
    "D42" = '1.001'
    "E42" = '  -0.04%  '
    "D43" = '10.50'
    "E43" = '  -4.61%  '
    "D44" = '0.5774'
    "E44" = '  -4.85%  '
    "D45" = '13.06'
    "E45" = '  -3.77%  '
    "D46" = '3.694'
    "E46" = '  -1.18%  '
    "D47" = '0.5528'
    "E47" = '  -6.19%  '
    "D48" = '1.887'
    "E48" = '  -5.23%  '
    "D49" = '117.22'
    "E49" = '  -5.07%  '
    "D50" = '1.133'
    "E50" = '  -3.55%  '
    "D51" = '0.06831'
    "E51" = '  -4.31%  '
}

function Test-LooksNumeric([string]$text) {
    # Matches optional sign, digits, optional single decimal point and more digits
    # (i.e. the shape of a plain decimal number with no thousands separators),
    # which is the shape Excel will silently auto-convert to a Number/Double.
    return $text -match '^[+-]?\d+(\.\d+)?$'
}

foreach ($ref in $newValues.Keys) {
    $value = $newValues[$ref]
    $range = $ws.Range($ref)

    if (Test-LooksNumeric $value) {
        # Force text storage so the cell keeps e.g. "1.003" rather than
        # becoming the number 1.003, then strip the quote-prefix style
        # Excel adds so the cell's formatting stays untouched.
        $range.Value = "'" + $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}
